$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.535604
$ws.Range("H2").Value = 49.606812
$ws.Range("I2").Value = 0.2120453146491552
$ws.Range("J2").Value = 0.2120453146491552
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.676256
$ws.Range("N2").Value = 5.028768
$ws.Range("O2").Value = 0.01781408474712955
$ws.Range("P2").Value = 0.01781408474712955
$ws.Range("Q2").Value = 27.717905418624
$ws.Range("R2").Value = 249.461148767616
$ws.Range("S2").Value = 0.003777393205391803
$ws.Range("T2").Value = 0.003777393205391803

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.535604
$ws.Range("H3").Value = 49.606812
$ws.Range("I3").Value = 0.2120453146491552
$ws.Range("J3").Value = 0.2120453146491552
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 63.45677799999999
$ws.Range("N3").Value = 190.370334
$ws.Range("O3").Value = 0.6743745711107288
$ws.Range("P3").Value = 0.6743745711107288
$ws.Range("Q3").Value = 1049.296152123912
$ws.Range("R3").Value = 9443.665369115206
$ws.Range("S3").Value = 0.1429979681225636
$ws.Range("T3").Value = 0.1429979681225636

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.535604
$ws.Range("H4").Value = 49.606812
$ws.Range("I4").Value = 0.2120453146491552
$ws.Range("J4").Value = 0.2120453146491552
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.96419433333334
$ws.Range("N4").Value = 86.892583
$ws.Range("O4").Value = 0.3078113441421415
$ws.Range("P4").Value = 0.3078113441421415
$ws.Range("Q4").Value = 478.940447675044
$ws.Range("R4").Value = 4310.464029075396
$ws.Range("S4").Value = 0.06526995332119979
$ws.Range("T4").Value = 0.06526995332119979

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.62063066666667
$ws.Range("H5").Value = 121.861892
$ws.Range("I5").Value = 0.5209011059384622
$ws.Range("J5").Value = 0.5209011059384622
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.676256
$ws.Range("N5").Value = 5.028768
$ws.Range("O5").Value = 0.01781408474712955
$ws.Range("P5").Value = 0.01781408474712955
$ws.Range("Q5").Value = 68.09057587878401
$ws.Range("R5").Value = 612.8151829090561
$ws.Range("S5").Value = 0.009279376446061275
$ws.Range("T5").Value = 0.009279376446061274

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 40.62063066666667
$ws.Range("H6").Value = 121.861892
$ws.Range("I6").Value = 0.5209011059384622
$ws.Range("J6").Value = 0.5209011059384622
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 63.45677799999999
$ws.Range("N6").Value = 190.370334
$ws.Range("O6").Value = 0.6743745711107288
$ws.Range("P6").Value = 0.6743745711107288
$ws.Range("Q6").Value = 2577.654342434659
$ws.Range("R6").Value = 23198.88908191193
$ws.Range("S6").Value = 0.3512824599083548
$ws.Range("T6").Value = 0.3512824599083548

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 40.62063066666667
$ws.Range("H7").Value = 121.861892
$ws.Range("I7").Value = 0.5209011059384622
$ws.Range("J7").Value = 0.5209011059384622
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 28.96419433333334
$ws.Range("N7").Value = 86.892583
$ws.Range("O7").Value = 0.3078113441421415
$ws.Range("P7").Value = 0.3078113441421415
$ws.Range("Q7").Value = 1176.543840571893
$ws.Range("R7").Value = 10588.89456514704
$ws.Range("S7").Value = 0.1603392695840461
$ws.Range("T7").Value = 0.1603392695840461

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 20.825229
$ws.Range("H8").Value = 62.475687
$ws.Range("I8").Value = 0.2670535794123827
$ws.Range("J8").Value = 0.2670535794123827
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.676256
$ws.Range("N8").Value = 5.028768
$ws.Range("O8").Value = 0.01781408474712955
$ws.Range("P8").Value = 0.01781408474712955
$ws.Range("Q8").Value = 34.908415062624
$ws.Range("R8").Value = 314.175735563616
$ws.Range("S8").Value = 0.004757315095676477
$ws.Range("T8").Value = 0.004757315095676476

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 20.825229
$ws.Range("H9").Value = 62.475687
$ws.Range("I9").Value = 0.2670535794123827
$ws.Range("J9").Value = 0.2670535794123827
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 63.45677799999999
$ws.Range("N9").Value = 190.370334
$ws.Range("O9").Value = 0.6743745711107288
$ws.Range("P9").Value = 0.6743745711107288
$ws.Range("Q9").Value = 1321.501933452162
$ws.Range("R9").Value = 11893.51740106946
$ws.Range("S9").Value = 0.1800941430798105
$ws.Range("T9").Value = 0.1800941430798105

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 20.825229
$ws.Range("H10").Value = 62.475687
$ws.Range("I10").Value = 0.2670535794123827
$ws.Range("J10").Value = 0.2670535794123827
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 28.96419433333334
$ws.Range("N10").Value = 86.892583
$ws.Range("O10").Value = 0.3078113441421415
$ws.Range("P10").Value = 0.3078113441421415
$ws.Range("Q10").Value = 603.1859797921691
$ws.Range("R10").Value = 5428.673818129521
$ws.Range("S10").Value = 0.08220212123689563
$ws.Range("T10").Value = 0.08220212123689562
